$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell. For cells whose new
# text happens to look like a plain number (e.g. "250.46"), Excel would
# otherwise auto-convert the literal string into a floating point number
# on assignment -- so those cells are switched to Text format first, and
# the format is reset back to the default "Normal" style afterwards so
# the cell keeps its original (unstyled) appearance.
function Set-LiteralText($range, [string]$text, [bool]$looksNumeric) {
    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-LiteralText $ws.Range("D2") "42.643.78" $false
Set-LiteralText $ws.Range("E2") "  -0.37%  " $false
Set-LiteralText $ws.Range("D3") "2.268.51" $false
Set-LiteralText $ws.Range("E3") "  +0.40%  " $false
Set-LiteralText $ws.Range("E4") "  -0.05%  " $false
Set-LiteralText $ws.Range("D5") "250.46" $true
Set-LiteralText $ws.Range("E5") "  +0.40%  " $false
Set-LiteralText $ws.Range("D6") "0.642" $true
Set-LiteralText $ws.Range("E6") "  +1.01%  " $false
Set-LiteralText $ws.Range("D7") "75.87" $true
Set-LiteralText $ws.Range("E7") "  +5.92%  " $false
Set-LiteralText $ws.Range("E8") "  -0.03%  " $false
Set-LiteralText $ws.Range("D9") "0.641" $true
Set-LiteralText $ws.Range("E9") "  -2.61%  " $false
Set-LiteralText $ws.Range("D10") "40.35" $true
Set-LiteralText $ws.Range("E10") "  +3.58%  " $false
Set-LiteralText $ws.Range("D11") "0.0973" $true
Set-LiteralText $ws.Range("E11") "  +0.29%  " $false
Set-LiteralText $ws.Range("E12") "  -1.31%  " $false
Set-LiteralText $ws.Range("E13") "  +0.93%  " $false
Set-LiteralText $ws.Range("D14") "2.611.70" $false
Set-LiteralText $ws.Range("E14") "  +0.76%  " $false
Set-LiteralText $ws.Range("D15") "15.03" $true
Set-LiteralText $ws.Range("E15") "  +0.85%  " $false
Set-LiteralText $ws.Range("D16") "0.866" $true
Set-LiteralText $ws.Range("E16") "  -1.65%  " $false
Set-LiteralText $ws.Range("D17") "2.281.18" $false
Set-LiteralText $ws.Range("E17") "  +0.72%  " $false
Set-LiteralText $ws.Range("D18") "42.550.73" $false
Set-LiteralText $ws.Range("E18") "  -0.43%  " $false
Set-LiteralText $ws.Range("D19") "0.0₃0993" $false
Set-LiteralText $ws.Range("E19") "  +0.01%  " $false
Set-LiteralText $ws.Range("D20") "6.18" $true
Set-LiteralText $ws.Range("E20") "  -2.15%  " $false
Set-LiteralText $ws.Range("D21") "72.10" $true
Set-LiteralText $ws.Range("E21") "  -1.39%  " $false
Set-LiteralText $ws.Range("D22") "235.94" $true
Set-LiteralText $ws.Range("E22") "  +0.95%  " $false
Set-LiteralText $ws.Range("D23") "2.16" $true
Set-LiteralText $ws.Range("E23") "  +3.37%  " $false
Set-LiteralText $ws.Range("E24") "  -2.22%  " $false
Set-LiteralText $ws.Range("E25") "  -0.01%  " $false
Set-LiteralText $ws.Range("D26") "11.19" $true
Set-LiteralText $ws.Range("E26") "  -2.56%  " $false
Set-LiteralText $ws.Range("E27") "  -1.82%  " $false
Set-LiteralText $ws.Range("D28") "2.20" $true
Set-LiteralText $ws.Range("E28") "  +4.17%  " $false
Set-LiteralText $ws.Range("D29") "167.41" $true
Set-LiteralText $ws.Range("E29") "  -0.23%  " $false
Set-LiteralText $ws.Range("D30") "20.91" $true
Set-LiteralText $ws.Range("E30") "  -0.23%  " $false
Set-LiteralText $ws.Range("D31") "6.45" $true
Set-LiteralText $ws.Range("E31") "  -1.80%  " $false
Set-LiteralText $ws.Range("D32") "0.0852" $true
Set-LiteralText $ws.Range("E32") "  +6.83%  " $false
Set-LiteralText $ws.Range("D33") "0.125" $true
Set-LiteralText $ws.Range("E33") "  -0.92%  " $false
Set-LiteralText $ws.Range("D34") "31.82" $true
Set-LiteralText $ws.Range("E34") "  +1.68%  " $false
Set-LiteralText $ws.Range("D35") "0.127" $true
Set-LiteralText $ws.Range("E35") "  +2.06%  " $false
Set-LiteralText $ws.Range("D36") "4.53" $true
Set-LiteralText $ws.Range("E36") "  +1.67%  " $false
Set-LiteralText $ws.Range("D37") "4.73" $true
Set-LiteralText $ws.Range("E37") "  -0.08%  " $false
Set-LiteralText $ws.Range("D38") "0.0305" $true
Set-LiteralText $ws.Range("E38") "  -4.70%  " $false
Set-LiteralText $ws.Range("D39") "13.43" $true
Set-LiteralText $ws.Range("E39") "  +6.79%  " $false
Set-LiteralText $ws.Range("D40") "2.27" $true
Set-LiteralText $ws.Range("E40") "  -2.02%  " $false
Set-LiteralText $ws.Range("D41") "5.86" $true
Set-LiteralText $ws.Range("E41") "  +0.70%  " $false
Set-LiteralText $ws.Range("D42") "0.206" $true
Set-LiteralText $ws.Range("E42") "  +1.89%  " $false
Set-LiteralText $ws.Range("D43") "61.32" $true
Set-LiteralText $ws.Range("E43") "  -1.40%  " $false
Set-LiteralText $ws.Range("D44") "8.89" $true
Set-LiteralText $ws.Range("E44") "  -2.89%  " $false
Set-LiteralText $ws.Range("D45") "106.13" $true
Set-LiteralText $ws.Range("E45") "  +11.60%  " $false
Set-LiteralText $ws.Range("D46") "4.72" $true
Set-LiteralText $ws.Range("E46") "  -2.70%  " $false
Set-LiteralText $ws.Range("D47") "0.101" $true
Set-LiteralText $ws.Range("E47") "  -1.54%  " $false
Set-LiteralText $ws.Range("E48") "  -0.52%  " $false
Set-LiteralText $ws.Range("D49") "1.15" $true
Set-LiteralText $ws.Range("E49") "  -0.96%  " $false
Set-LiteralText $ws.Range("E50") "  -2.35%  " $false
Set-LiteralText $ws.Range("D51") "4.18" $true
Set-LiteralText $ws.Range("E51") "  -2.94%  " $false
